$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Update the IP address for the ZVIU-D002/CDLU-D001 stream encoder
# (row 46) -- it now has a known reserved IP instead of "tbd".
# -----------------------------------------------------------------
$ws.Range("I46").Value = "192.168.0.103"

# -----------------------------------------------------------------
# New equipment: Verkada door-access controller, added to the
# bottom of the table as rows 52-53.
# -----------------------------------------------------------------

# Row 52 -- asset-tag / category / MAC-only placeholder line
$ws.Range("A52").Value = "ZVIU-G001"
$ws.Range("A52").NumberFormat = "@"
$ws.Range("B52").Value = "Video"

$ws.Range("G52").Value = "1C-69-7A-66-87-69"
$g52Font = $ws.Range("G52").Font
$g52Font.Size = 9
$g52Font.Name = "Helvetica"

# Row 53 -- the actual device row, with a top/bottom separator
# border to set it apart from the rest of the table
$ws.Range("A53").Value = "2405-1307"
$ws.Range("A53").NumberFormat = "@"

$topBorder = $ws.Range("A53").Borders.Item(8)
$topBorder.Color = 14461583
$topBorder.Weight = 2
$topBorder.LineStyle = 1

$bottomBorder = $ws.Range("A53").Borders.Item(9)
$bottomBorder.Color = 14461583
$bottomBorder.Weight = 2
$bottomBorder.LineStyle = 1

$ws.Range("B53").Value = "Security"
$ws.Range("E53").Value = "Door Access Controller"
$ws.Range("F53").Value = "Verkada Door Control"
$ws.Range("G53").Value = "E0:A7:00:3F:55:77"
$ws.Range("I53").Value = "192.168.0.239"
$ws.Range("K53").Value = "amp room"

# Move the cursor/selection onto the new last row, matching where
# the author left off editing.
$ws.Range("K53").Select() | Out-Null

Write-Output "done"
